# TestComposeKey.xlsx — "Add DefaultAttr & AliasAttr"
#
# The "Pos" field's type-annotation cell (D2) gains a second line declaring
# an alias attribute, e.g. the row layout is:
#   Row1 (field name)  : ... | Pos           | ...
#   Row2 (type/attrs)  : ... | float_float_float  -> float_float_float / Alias[V3]
#   Row3 (display name): ... | 位置           | ...
#
# so D2's text becomes two lines: "float_float_float" + "Alias[V3]".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Core content edit: append the new Alias[V3] attribute line to D2.
$ws.Range("D2").Value = "float_float_float`r`nAlias[V3]"

# The new attribute text is noticeably longer, so column D is widened to
# fit it (matches the new explicit column width in the sheet).
$ws.Columns.Item(4).ColumnWidth = 25.75

# Reflect the author's last selection/cursor position on the sheet.
$ws.Range("K6").Select()
